$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.81135368347168
$ws.Range("B1").Value = 2.574819326400757
$ws.Range("C1").Value = 1.845252394676208
$ws.Range("D1").Value = 1.687134265899658
$ws.Range("E1").Value = 1.706208944320679
